# Add 2022-Q1 data:
#  - the existing "总计" sheet becomes "2022-Q1" (new fund-holdings data)
#  - a brand-new "总计" sheet is appended with the updated totals history

$wb = $excel.ActiveWorkbook

$q1sheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet as the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1sheet.Name = "2022-Q1"
$q1sheet.Cells.Clear()

# Bring over the header-row / index-column formatting from the most
# recent quarter sheet (2021-Q4) so the new sheet matches the
# established visual style (bold header, bordered index column, etc).
$templateSheet.Range("B1:H1").Copy() | Out-Null
$q1sheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$templateSheet.Range("A2:A6").Copy() | Out-Null
$q1sheet.Range("A2:A8").PasteSpecial(-4122) | Out-Null

$q1sheet.Range("B1").Value = "基金代码"
$q1sheet.Range("C1").Value = "基金名称"
$q1sheet.Range("D1").Value = "基金规模"
$q1sheet.Range("E1").Value = "股票总仓位"
$q1sheet.Range("F1").Value = "仓位占比"
$q1sheet.Range("G1").Value = "持有市值(亿元)"
$q1sheet.Range("H1").Value = "仓位排名"

# Columns B-G hold text (fund codes must keep leading zeros, and the
# numeric-looking figures must stay literal text, matching the source).
$q1sheet.Range("B2:G8").NumberFormat = "@"

$q1data = @(
    @(0, "516950", "银华中证基建交易型开放式指数证券投资基金", "10.41", "97.55", "3.62", "0.3768", 10),
    @(1, "161123", "易方达并购重组指数（LOF）", "4.78", "94.71", "4.57", "0.2184", 5),
    @(2, "166802", "浙商沪深 300 指数增强（LOF）", "4.25", "88.42", "1.88", "0.0799", 10),
    @(3, "002025", "广发聚盛灵活配置混合A", "7.09", "22.40", "0.29", "0.0206", 10),
    @(4, "515870", "嘉实中证先进制造100策略ETF", "0.42", "98.79", "3.25", "0.0136", 10),
    @(5, "005502", "华泰紫金智能量化股票", "0.43", "94.38", "2.02", "0.0087", 5),
    @(6, "002026", "广发聚盛灵活配置混合C", "1.07", "22.40", "0.29", "0.0031", 10)
)

$r = 2
foreach ($row in $q1data) {
    $q1sheet.Cells.Item($r, 1).Value = $row[0]
    $q1sheet.Cells.Item($r, 2).Value = $row[1]
    $q1sheet.Cells.Item($r, 3).Value = $row[2]
    $q1sheet.Cells.Item($r, 4).Value = $row[3]
    $q1sheet.Cells.Item($r, 5).Value = $row[4]
    $q1sheet.Cells.Item($r, 6).Value = $row[5]
    $q1sheet.Cells.Item($r, 7).Value = $row[6]
    $q1sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet with the updated totals history
# ---------------------------------------------------------------------
# Copy a whole sheet (rather than Worksheets.Add()) so sheet-level
# properties (sheetPr/outline, page margins, etc.) match the rest of
# the workbook instead of falling back to blank-sheet defaults.
$templateSheet.Copy($null, $q1sheet) | Out-Null
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalSheet.Name = "总计"
$totalSheet.Cells.Clear()

$templateSheet.Range("B1:D1").Copy() | Out-Null
$totalSheet.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$templateSheet.Range("A2:A6").Copy() | Out-Null
$totalSheet.Range("A2:A7").PasteSpecial(-4122) | Out-Null

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @(0, "2022-Q1", 7, 0.72),
    @(1, "2021-Q4", 5, 0.9),
    @(2, "2021-Q3", 22, 4.49),
    @(3, "2021-Q2", 7, 3.55),
    @(4, "2021-Q1", 8, 4.25),
    @(5, "2020-Q4", 5, 0.23)
)

$r = 2
foreach ($row in $totalData) {
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Restore the original active sheet (the copy/activation above moves
# Excel's focus to the newly created sheets).
$wb.Worksheets.Item(1).Activate()
